$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.167.43"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "1.828.98"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.14"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6200"
$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07361"
$ws.Range("E8").Value = "  -1.90%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2905"
$ws.Range("E9").Value = "  -1.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.09"
$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07682"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").Value = "1.822.46"
$ws.Range("E12").Value = "  -0.82%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.954"
$ws.Range("E13").Value = "  -1.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6641"
$ws.Range("E14").Value = "  -1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.18"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008911"
$ws.Range("E16").Value = "  -4.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.853"
$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("D18").Value = "29.136.02"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("D19").Value = "2.065.06"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.44"
$ws.Range("E20").Value = "  +5.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("E23").Value = "  +2.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.82"
$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("E26").Value = "  +0.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.494"
$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.65"
$ws.Range("E28").Value = "  -1.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.487"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05909"
$ws.Range("E30").Value = "  +6.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.073"
$ws.Range("E31").Value = "  -1.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.079"
$ws.Range("E32").Value = "  -2.48%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.861"
$ws.Range("E34").Value = "  +0.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7313"
$ws.Range("E35").Value = "  -2.38%  "

$ws.Range("E36").Value = "  -0.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.614"
$ws.Range("E37").Value = "  -1.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.846"
$ws.Range("E38").Value = "  +2.64%  "

$ws.Range("D39").Value = "1.216.80"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01745"
$ws.Range("E40").Value = "  -2.30%  "

$ws.Range("E41").Value = "  -4.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9159"
$ws.Range("E42").Value = "  +2.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.76"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("D45").Value = "1.968.29"
$ws.Range("E45").Value = "  -0.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.90"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5085"
$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.159"
$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4020"

$ws.Range("E50").Value = "  -4.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1128"
$ws.Range("E51").Value = "  +2.10%  "
